$d = $word.ActiveDocument

$old = "Informace v této příručce jsou určeny pro pozorovací kampaň probíhající od Herkules souhvězdí 2022: 13. – 22. června, 12. – 21. července, 10. – 19. srpna"
$new = "Informace v této příručce jsou určeny pro pozorovací kampaň probíhající od 13. – 22. června, 12. – 21. července, 10. – 19. srpna. Při pozorování použijte hvězdy oblohy, které zobrazují souhvězdí Herkules souhvězdí.13. – 22. června, 12. – 21. července, 10. – 19. srpna"

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
